$wb = $excel.ActiveWorkbook

# Add the new "calendar" worksheet after the last existing sheet ("company")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "calendar"

# Header row - reuse the same labels ("title"/"category") used on the other sheets
$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "category"
$ws.Range("A1:B1").Interior.Color = 65535

# Data row - new calendar/event test case
$ws.Range("A2").Value = "Test title - 1"
$ws.Range("B2").Value = "Important"

# Column widths approximating the source workbook's autofit sizing
$ws.Columns.Item(1).ColumnWidth = 10.6667
$ws.Columns.Item(2).ColumnWidth = 9

# Selection as recorded in the edited workbook
$ws.Range("C5").Select()
